# Atualização automática: 2025-08-19 21:00:25
#
# Rows 7-11 of the detections table are rotated by one position (the record
# that used to live on row 11 moves up to row 7, and the records that used to
# occupy rows 7-10 each shift down by one row). Rows 16-18 get refreshed
# detection-image / coordinate / confidence values for a later pass of the
# same image.
#
# Columns I (First_Coords, e.g. "962,713,1006,765") and J (First_Confidence,
# e.g. "0.76") are plain text in the source data, so we force the cell format
# to Text ("@") before writing them - otherwise Excel would either treat the
# comma separated coordinate quadruplet as a grouped number (stripping the
# commas) or convert the confidence score into a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (becomes the former row 11's record) ---
$ws.Range("A7").Value = "2117575c-4ae1-458c-b88a-fc40f40debdb"
$ws.Range("D7").Value = "image_20250727074723_ppp0.jpg"
$ws.Range("E7").Value = "PLACA_20250723145134"
$ws.Range("F7").Value = "Moura"
$ws.Range("G7").Value = 38.06587
$ws.Range("H7").Value = -7.221796
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "1490,161,1563,258"
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = "0.62"

# --- Row 8 (becomes the former row 7's record) ---
$ws.Range("A8").Value = "283b6eda-9c83-4cdd-9524-c7c394f2dc89"
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "962,713,1006,765"
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "0.76"

# --- Row 9 (becomes the former row 8's record) ---
$ws.Range("A9").Value = "a19b65d1-6f97-4841-9e1c-7446a9be92b6"
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "967,614,1002,659"
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "0.73"

# --- Row 10 (becomes the former row 9's record) ---
$ws.Range("A10").Value = "4be1b1cf-d480-453e-b5fb-d4ecd6764c4d"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "702,633,740,690"
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = "0.72"

# --- Row 11 (becomes the former row 10's record) ---
$ws.Range("A11").Value = "dfd476d4-7689-4671-a076-78fe3ce806bb"
$ws.Range("D11").Value = "image_20250728214139_ppp0.jpg"
$ws.Range("E11").Value = "PLACA_20250717165933"
$ws.Range("F11").Value = "Beja"
$ws.Range("G11").Value = 38.02035
$ws.Range("H11").Value = -7.94715
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "1254,850,1294,895"
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = "0.67"

# --- Row 16: refreshed detection image / bbox / confidence ---
$ws.Range("D16").Value = "image_20250807111314_ppp0.jpg"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "643,531,686,575"
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "0.76"

# --- Row 17: refreshed detection image / bbox / confidence ---
$ws.Range("D17").Value = "image_20250807111314_ppp0.jpg"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "794,481,830,526"
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = "0.72"

# --- Row 18: refreshed detection image / bbox / confidence ---
$ws.Range("D18").Value = "image_20250808100711_ppp0.jpg"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "1182,409,1232,451"
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = "0.75"
